$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Septiembre de 2020 a las 21:30"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 7167280
$ws.Range("C4").Value = 27727
$ws.Range("D4").Value = 4415825
$ws.Range("E4").Value = 2544274
$ws.Range("G4").Value = 588
$ws.Range("H4").Value = 207181

# --- Row 5: India ---
$ws.Range("B5").Value = 5813799
$ws.Range("C5").Value = 83615
$ws.Range("D5").Value = 4750922
$ws.Range("E5").Value = 970584
$ws.Range("G5").Value = 1120
$ws.Range("H5").Value = 92293

# --- Row 6: Brasil ---
$ws.Range("B6").Value = 4634468
$ws.Range("C6").Value = 6688
$ws.Range("E6").Value = 502288
$ws.Range("G6").Value = 229
$ws.Range("H6").Value = 139294

# --- Row 12: Sudafrica ---
$ws.Range("B12").Value = 667049
$ws.Range("C12").Value = 1861
$ws.Range("D12").Value = 595916
$ws.Range("E12").Value = 54850
$ws.Range("G12").Value = 77
$ws.Range("H12").Value = 16283

# --- Row 14: Francia ---
$ws.Range("B14").Value = 497237
$ws.Range("C14").Value = 16096
$ws.Range("D14").Value = 94413
$ws.Range("E14").Value = 371313
$ws.Range("G14").Value = 52
$ws.Range("H14").Value = 31511

# --- Row 25: Alemania ---
$ws.Range("B25").Value = 280853
$ws.Range("C25").Value = 1648
$ws.Range("E25").Value = 21836
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 9517

# --- Row 31: Ecuador ---
$ws.Range("B31").Value = 131146
$ws.Range("C31").Value = 1254
$ws.Range("E31").Value = 17081
$ws.Range("G31").Value = 42
$ws.Range("H31").Value = 11213

# --- Row 75: Libano ---
$ws.Range("B75").Value = 32819
$ws.Range("C75").Value = 1027
$ws.Range("D75").Value = 14112
$ws.Range("E75").Value = 18378
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 329

# --- Row 101: Guayana Francesa ---
$ws.Range("B101").Value = 9790
$ws.Range("C101").Value = 28
$ws.Range("D101").Value = 9456
$ws.Range("E101").Value = 269

# --- Row 117: Cabo Verde ---
$ws.Range("B117").Value = 5479
$ws.Range("C117").Value = 67
$ws.Range("D117").Value = 4917
$ws.Range("E117").Value = 507
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 55

# --- Row 144: Mali ---
$ws.Range("B144").Value = 3041
$ws.Range("C144").Value = 7
$ws.Range("D144").Value = 2391
$ws.Range("E144").Value = 520

# --- Row 190: Monaco ---
$ws.Range("B190").Value = 205
$ws.Range("C190").Value = 6
$ws.Range("D190").Value = 166
$ws.Range("E190").Value = 38

# --- Rows 215/216: swap Islas Malvinas <-> Montserrat (name + their D/H stats) ---
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
